$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 5133205
$ws.Range("I43").Value = 7694807.5
$ws.Range("K43").Value = 7694807.5
$ws.Range("M43").Value = -7694738.5
$ws.Range("H55").Value = 654.8
$ws.Range("I55").Value = 623.875
$ws.Range("J55").Value = 778.5
$ws.Range("K55").Value = 623.875
$ws.Range("L55").Value = 778.5
$ws.Range("M55").Value = -409.875
$ws.Range("N55").Value = -1206.5
$ws.Range("H100").Value = 1067
$ws.Range("I100").Value = 1289.4
$ws.Range("J100").Value = 696.3333
$ws.Range("K100").Value = 1289.4
$ws.Range("L100").Value = 696.3333
$ws.Range("M100").Value = -748.4000000000001
$ws.Range("N100").Value = -1778.3333
$ws.Range("H121").Value = 4047.8965
$ws.Range("J121").Value = 4047.8965
$ws.Range("L121").Value = 12143.6895
$ws.Range("N121").Value = -15637.6895
$ws.Range("H138").Value = 4761
$ws.Range("I138").Value = 1436
$ws.Range("J138").Value = 5276.948
$ws.Range("K138").Value = 4308
$ws.Range("L138").Value = 15830.844
$ws.Range("M138").Value = 832
$ws.Range("N138").Value = -26110.844

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5363.4424
$ws.Range("I32").Value = 3006.4634
$ws.Range("K32").Value = 3006.4634
$ws.Range("M32").Value = -2719.4634
$ws.Range("H45").Value = 53455.137
$ws.Range("I45").Value = 67707.12
$ws.Range("K45").Value = 67707.12
$ws.Range("M45").Value = -67330.12
$ws.Range("H132").Value = 5014.6855
$ws.Range("I132").Value = 1695.5264
$ws.Range("J132").Value = 8956.1875
$ws.Range("K132").Value = 5086.5792
$ws.Range("L132").Value = 26868.5625
$ws.Range("M132").Value = -2556.5792
$ws.Range("N132").Value = -31928.5625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1359.238
$ws.Range("I86").Value = 1376.5333
$ws.Range("J86").Value = 1316
$ws.Range("K86").Value = 1376.5333
$ws.Range("L86").Value = 1316
$ws.Range("M86").Value = -253.5333000000001
$ws.Range("N86").Value = -3562
$ws.Range("H89").Value = 1359.238
$ws.Range("I89").Value = 1376.5333
$ws.Range("J89").Value = 1316
$ws.Range("K89").Value = 6882.6665
$ws.Range("L89").Value = 6580
$ws.Range("M89").Value = -1266.6665
$ws.Range("N89").Value = -17812
$ws.Range("H134").Value = 1042.238
$ws.Range("I134").Value = 983.8421
$ws.Range("K134").Value = 2951.5263
$ws.Range("M134").Value = -416.5263

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2444.3333
$ws.Range("J105").Value = 2444.3333
$ws.Range("L105").Value = 2444.3333
$ws.Range("N105").Value = -5938.3333
$ws.Range("H132").Value = 17561222
$ws.Range("I132").Value = 18531290
$ws.Range("J132").Value = 100000
$ws.Range("K132").Value = 55593870
$ws.Range("L132").Value = 300000
$ws.Range("M132").Value = -55591340
$ws.Range("N132").Value = -305060
$ws.Range("H134").Value = 2057.8
$ws.Range("I134").Value = 2066.6155
$ws.Range("K134").Value = 6199.8465
$ws.Range("M134").Value = -3664.8465
$ws.Range("H141").Value = 95019.64
$ws.Range("J141").Value = 98998.16
$ws.Range("L141").Value = 98998.16
$ws.Range("N141").Value = -109358.16

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2920671.8
$ws.Range("I4").Value = 1125443.2
$ws.Range("K4").Value = 3376329.6
$ws.Range("M4").Value = -3376217.6
$ws.Range("H39").Value = 3148.3333
$ws.Range("J39").Value = 2250
$ws.Range("L39").Value = 6750
$ws.Range("N39").Value = -7338
$ws.Range("H56").Value = 6649.8887
$ws.Range("I56").Value = 6649.8887
$ws.Range("K56").Value = 6649.8887
$ws.Range("M56").Value = -6119.8887
$ws.Range("H86").Value = 141.53847
$ws.Range("J86").Value = 100
$ws.Range("L86").Value = 300
$ws.Range("N86").Value = -2672
$ws.Range("H89").Value = 141.53847
$ws.Range("J89").Value = 100
$ws.Range("L89").Value = 900
$ws.Range("N89").Value = -12756
$ws.Range("H140").Value = 4009
$ws.Range("I140").Value = 2348.2778
$ws.Range("K140").Value = 7044.8334
$ws.Range("M140").Value = -1864.8334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 60360.684
$ws.Range("I80").Value = 128023.875
$ws.Range("J80").Value = 21696
$ws.Range("K80").Value = 128023.875
$ws.Range("L80").Value = 21696
$ws.Range("M80").Value = -127025.875
$ws.Range("N80").Value = -23692
$ws.Range("H83").Value = 60360.684
$ws.Range("I83").Value = 128023.875
$ws.Range("J83").Value = 21696
$ws.Range("K83").Value = 640119.375
$ws.Range("L83").Value = 108480
$ws.Range("M83").Value = -635127.375
$ws.Range("N83").Value = -118464
$ws.Range("H97").Value = 495.81818
$ws.Range("I97").Value = 484.25
$ws.Range("J97").Value = 526.6667
$ws.Range("K97").Value = 484.25
$ws.Range("L97").Value = 526.6667
$ws.Range("M97").Value = 11.75
$ws.Range("N97").Value = -1518.6667
$ws.Range("H126").Value = 5864.273
$ws.Range("I126").Value = 3748.75
$ws.Range("K126").Value = 11246.25
$ws.Range("M126").Value = -8776.25
$ws.Range("H132").Value = 206200.1
$ws.Range("I132").Value = 339531.34
$ws.Range("J132").Value = 6203.25
$ws.Range("K132").Value = 1018594.02
$ws.Range("L132").Value = 18609.75
$ws.Range("M132").Value = -1016064.02
$ws.Range("N132").Value = -23669.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6292.3447
$ws.Range("J46").Value = 6388.074
$ws.Range("L46").Value = 6388.074
$ws.Range("N46").Value = -6764.074
$ws.Range("H122").Value = 51953504
$ws.Range("I122").Value = 125004250
$ws.Range("K122").Value = 375012750
$ws.Range("M122").Value = -375010300
$ws.Range("H132").Value = 4208.0356
$ws.Range("I132").Value = 3645.5652
$ws.Range("K132").Value = 10936.6956
$ws.Range("M132").Value = -8406.6956
$ws.Range("H136").Value = 5884.615
$ws.Range("I136").Value = 2400
$ws.Range("J136").Value = 8062.5
$ws.Range("K136").Value = 7200
$ws.Range("L136").Value = 24187.5
$ws.Range("M136").Value = -4650
$ws.Range("N136").Value = -29287.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 22227566
$ws.Range("I132").Value = 3969037
$ws.Range("K132").Value = 11907111
$ws.Range("M132").Value = -11904581
$ws.Range("H139").Value = 59785.75
$ws.Range("J139").Value = 59785.75
$ws.Range("L139").Value = 59785.75
$ws.Range("N139").Value = -70065.75

Write-Output "Applied all edits"